$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# Edit 1: "Can we allow the user to suggest the location and help to display
# the location." -> split the run containing "location" (2nd occurrence) off
# into its own run wrapped in proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Can we allow the user to suggest the location and help to display the location.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: paragraph 1 (Can we allow...) not found"
}
$para = $rng.Paragraphs(1)
$prng = $para.Range

$rPr1 = '<w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="000000"/></w:rPr>'
$pPr1 = '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="280" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/><w:color w:val="000000"/></w:rPr></w:pPr>'

$p1 = '<w:p ' + $wns + ' ' + $w14ns + ' w14:paraId="6D04F2C4" w14:textId="1F16576F" w:rsidR="00A87065" w:rsidRDefault="00A87065" w:rsidP="00CC6514">'
$p1 += $pPr1
$p1 += '<w:r>' + $rPr1 + '<w:t>Can we allow the user to sug</w:t></w:r>'
$p1 += '<w:r>' + $rPr1 + '<w:t xml:space="preserve">gest the location and help to display the </w:t></w:r>'
$p1 += '<w:proofErr w:type="gramStart"/>'
$p1 += '<w:r>' + $rPr1 + '<w:t>location</w:t></w:r>'
$p1 += '<w:r>' + $rPr1 + '<w:t>.</w:t></w:r>'
$p1 += '<w:proofErr w:type="gramEnd"/>'
$p1 += '</w:p>'

$prng.InsertXML($p1)
Write-Host "Edit 1 done"

# ---------------------------------------------------------------------------
# Edit 2: remove the _GoBack bookmark from the "Are we validating through
# NRIC?" paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Are we validating through NRIC?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: paragraph 2 (Are we validating...) not found"
}
$para = $rng.Paragraphs(1)
$prng = $para.Range

$p2 = '<w:p ' + $wns + ' ' + $w14ns + ' w14:paraId="204CED28" w14:textId="34C1F972" w:rsidR="00F1786B" w:rsidRDefault="00F1786B">'
$p2 += '<w:r><w:t>Are we validating through NRIC?</w:t></w:r>'
$p2 += '</w:p>'

$prng.InsertXML($p2)
Write-Host "Edit 2 done"

# ---------------------------------------------------------------------------
# Edit 3: "Wanying -> Image recognition and barcode scanning -> Web
# scraping" -> wrap "Wanying" with proofErr spellStart/spellEnd and split
# into two runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Wanying -> Image recognition and barcode scanning -> Web scraping", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: paragraph 3 (Wanying...) not found"
}
$para = $rng.Paragraphs(1)
$prng = $para.Range

$p3 = '<w:p ' + $wns + ' ' + $w14ns + ' w14:paraId="363324E9" w14:textId="631BCE97" w:rsidR="00392AFB" w:rsidRDefault="00392AFB">'
$p3 += '<w:proofErr w:type="spellStart"/>'
$p3 += '<w:r><w:t>Wanying</w:t></w:r>'
$p3 += '<w:proofErr w:type="spellEnd"/>'
$p3 += '<w:r><w:t xml:space="preserve"> -&gt; Image recognition and barcode scanning -&gt; Web scraping</w:t></w:r>'
$p3 += '</w:p>'

$prng.InsertXML($p3)
Write-Host "Edit 3 done"

# ---------------------------------------------------------------------------
# Edit 4: append a large block of new paragraphs after "Suggest UI, UX
# themes" (the last paragraph of the body), ending with the _GoBack
# bookmark that used to sit elsewhere in the document.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Suggest UI, UX themes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: paragraph 4 (Suggest UI, UX themes) not found"
}
$para = $rng.Paragraphs(1)
$prng = $para.Range

$lang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

$big = '<w:p ' + $wns + ' ' + $w14ns + ' w14:paraId="5EC9A482" w14:textId="78B905D6" w:rsidR="009C434C" w:rsidRDefault="008F097C">'
$big += '<w:r><w:t>Suggest UI, UX themes</w:t></w:r>'
$big += '</w:p>'

$big += '<w:p/>'

$big += '<w:p><w:pPr><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr>'
$big += '<w:r><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>26</w:t></w:r>'
$big += '<w:r><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>'
$big += '<w:r><w:rPr><w:b/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> April 2017</w:t></w:r>'
$big += '</w:p>'

$big += '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Can we also do filtering in the school, donating books back to the school and whatnot</w:t></w:r></w:p>'

$big += '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r>' + $lang + '<w:t>Personal and corporate</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Navigate when meeting is soon -&gt; open auto navigation with manual entry or current location?</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t xml:space="preserve">Data filtering based on personal preference (Prefer Sci-Fi/Fantasy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r>' + $lang + '<w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>' + $lang + '<w:t>)</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>E-mail</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Password</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>User name/Name</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>With preference of using username/name</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>User profile -&gt; Key in multiple schools</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>E-mail verification</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Chat bot for information</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>How do you keep a sustained user base?</w:t></w:r></w:p>'

$big += '<w:p><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Think up of a UI theme</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Adobe XD</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Try to do barcode scanning basic first</w:t></w:r><w:r>' + $lang + '<w:t xml:space="preserve"> (portable code)</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Home page</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Login</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Registration</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Profile</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Donation</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Scanning</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Chat</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t xml:space="preserve">Client </w:t></w:r><w:r>' + $lang + '<w:t>C</w:t></w:r><w:r>' + $lang + '<w:t>onnect</w:t></w:r><w:r>' + $lang + '<w:t>ion</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Map Route</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Browsing</w:t></w:r></w:p>'
$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:r>' + $lang + '<w:t>Searching</w:t></w:r></w:p>'

$big += '<w:p><w:pPr>' + $lang + '</w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$prng.InsertXML($big)
Write-Host "Edit 4 done"

Write-Host "All document.xml edits applied"
